$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.116974472999573
$ws.Range("B1").Value = 1.931352734565735
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.113637447357178
$ws.Range("E1").Value = 1.209613561630249
